$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 3 (the "ZA7750 / 94.2" row),
# shifting it (and everything below) down by one.
$ws.Rows.Item(3).Insert()

# Populate the new row 3 with the new survey wave metadata.
$ws.Range("A3").Value() = "ZA7780"
$ws.Range("B3").Value() = "'94.3"
$ws.Range("C3").Value() = "February - March 2021"
$ws.Range("D3").Value() = "COVID-19 Pandemic"

# Match the author's final selection (cell D4, the description cell of the
# row that got pushed down from 3 to 4).
$ws.Range("D4").Select()
